$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 310 }

# New "Förändrad" (Changed) date serial number for all data rows.
$newDateSerial = 45186

# Columns that hold HYPERLINK() formulas needing a friendly-name 2nd argument,
# but only for the rows that actually contain such formulas (rows 2-11).
$hyperlinkCols = @("S", "T", "V", "W", "X", "Y")

for ($row = 2; $row -le $lastRow; $row++) {

    # Update column C (Förändrad) on every row that has a value there.
    $cCell = $ws.Range("C$row")
    if ($cCell.Value2 -ne $null) {
        $cCell.Value = $newDateSerial
    }

    # Only the first 10 data rows (2-11) carry the hyperlink formulas.
    if ($row -ge 2 -and $row -le 11) {
        $label = $ws.Range("A$row").Value2

        foreach ($col in $hyperlinkCols) {
            $cell = $ws.Range("$col$row")
            $formula = $cell.Formula
            if ([string]::IsNullOrEmpty($formula)) { continue }

            if ($formula -match '^\=HYPERLINK\(("[^"]*")\)$') {
                $urlArg = $Matches[1]
                $cell.Formula = '=HYPERLINK(' + $urlArg + ', "' + $label + '")'
            }
        }
    }
}
